# Update hotel reviews data (review_info / hotel_info sheets)
$wb = $excel.ActiveWorkbook

# --- hotel_info sheet: set English_Reviews_num (G2) to '14' as TEXT ---
$ws1 = $wb.Worksheets.Item(1)
$c = $ws1.Cells.Item(2,7)
$c.NumberFormat = "@"
$c.Value = '14'
$c.Style = "Normal"

# --- review_info sheet: append 4 new review rows ---
$ws2 = $wb.Worksheets.Item(2)

# row 2
$ws2.Cells.Item(2,1).Value = 39431
$ws2.Cells.Item(2,4).Value = 1
$c = $ws2.Cells.Item(2,5)
$c.NumberFormat = "@"
$c.Value = '08/03/2018'
$c.Style = "Normal"
$c = $ws2.Cells.Item(2,6)
$c.NumberFormat = "@"
$c.Value = 'https://www.tripadvisor.com/ShowUserReviews-g29144-d119933-r277920100-InTown_Suites_Denver_East-Aurora_Colorado.html'
$c.Style = "Normal"
$c = $ws2.Cells.Item(2,7)
$c.NumberFormat = "@"
$c.Value = '29144'
$c.Style = "Normal"
$c = $ws2.Cells.Item(2,8)
$c.NumberFormat = "@"
$c.Value = '119933'
$c.Style = "Normal"
$c = $ws2.Cells.Item(2,9)
$c.NumberFormat = "@"
$c.Value = '277920100'
$c.Style = "Normal"
$c = $ws2.Cells.Item(2,10)
$c.NumberFormat = "@"
$c.Value = '06/04/2015'
$c.Style = "Normal"
$c = $ws2.Cells.Item(2,11)
$c.NumberFormat = "@"
$c.Value = 'Very disturbing experience'
$c.Style = "Normal"
$c = $ws2.Cells.Item(2,12)
$c.NumberFormat = "@"
$c.Value = 'Do not stay here!  I repeat, do not stay here!   I literally witnessed the manager doing dope with some very unsavory characters.  We got no sleep at all because of all the noise from doors slamming and people yelling.  I could not recommend this place to anybody.  It was a very sad commentary on how things have declined in communities affected by addiction.'
$c.Style = "Normal"
$ws2.Cells.Item(2,13).Value = 1
$c = $ws2.Cells.Item(2,14)
$c.NumberFormat = "@"
$c.Value = 'May 2015'
$c.Style = "Normal"
$c = $ws2.Cells.Item(2,15)
$c.NumberFormat = "@"
$c.Value = ' traveled as a couple'
$c.Style = "Normal"
$ws2.Cells.Item(2,19).Value = 1
$ws2.Cells.Item(2,21).Value = 1
$ws2.Cells.Item(2,22).Value = 0
$c = $ws2.Cells.Item(2,25)
$c.NumberFormat = "@"
$c.Value = 'Do not stay here!  I repeat, do not stay here!   I literally witnessed the manager doing dope with some very unsavory characters.  We got no sleep at all because of all the noise from doors slamming and people yelling.  I could not recommend this place to anybody.  It was a very sad commentary on how things have declined in communities affected by addiction.'
$c.Style = "Normal"

# row 3
$ws2.Cells.Item(3,1).Value = 39431
$ws2.Cells.Item(3,4).Value = 2
$c = $ws2.Cells.Item(3,5)
$c.NumberFormat = "@"
$c.Value = '08/03/2018'
$c.Style = "Normal"
$c = $ws2.Cells.Item(3,6)
$c.NumberFormat = "@"
$c.Value = 'https://www.tripadvisor.com/ShowUserReviews-g29144-d119933-r269124661-InTown_Suites_Denver_East-Aurora_Colorado.html'
$c.Style = "Normal"
$c = $ws2.Cells.Item(3,7)
$c.NumberFormat = "@"
$c.Value = '29144'
$c.Style = "Normal"
$c = $ws2.Cells.Item(3,8)
$c.NumberFormat = "@"
$c.Value = '119933'
$c.Style = "Normal"
$c = $ws2.Cells.Item(3,9)
$c.NumberFormat = "@"
$c.Value = '269124661'
$c.Style = "Normal"
$c = $ws2.Cells.Item(3,10)
$c.NumberFormat = "@"
$c.Value = '04/30/2015'
$c.Style = "Normal"
$c = $ws2.Cells.Item(3,11)
$c.NumberFormat = "@"
$c.Value = 'Awful Management!'
$c.Style = "Normal"
$c = $ws2.Cells.Item(3,12)
$c.NumberFormat = "@"
$c.Value = 'Luke warm reception when booking the room and then it went totally down hill from there.  General Manager, Jesse Zuniga, yells at customers using swear words when they have valid concerns.  Place is full of drugs, thugs and bugs.'
$c.Style = "Normal"
$ws2.Cells.Item(3,13).Value = 1
$c = $ws2.Cells.Item(3,14)
$c.NumberFormat = "@"
$c.Value = 'April 2015'
$c.Style = "Normal"
$c = $ws2.Cells.Item(3,15)
$c.NumberFormat = "@"
$c.Value = ' traveled solo'
$c.Style = "Normal"
$ws2.Cells.Item(3,19).Value = 1
$ws2.Cells.Item(3,21).Value = 1
$ws2.Cells.Item(3,22).Value = 0
$c = $ws2.Cells.Item(3,25)
$c.NumberFormat = "@"
$c.Value = 'Luke warm reception when booking the room and then it went totally down hill from there.  General Manager, Jesse Zuniga, yells at customers using swear words when they have valid concerns.  Place is full of drugs, thugs and bugs.'
$c.Style = "Normal"

# row 4
$ws2.Cells.Item(4,1).Value = 39431
$ws2.Cells.Item(4,4).Value = 3
$c = $ws2.Cells.Item(4,5)
$c.NumberFormat = "@"
$c.Value = '08/03/2018'
$c.Style = "Normal"
$c = $ws2.Cells.Item(4,6)
$c.NumberFormat = "@"
$c.Value = 'https://www.tripadvisor.com/ShowUserReviews-g29144-d119933-r260194944-InTown_Suites_Denver_East-Aurora_Colorado.html'
$c.Style = "Normal"
$c = $ws2.Cells.Item(4,7)
$c.NumberFormat = "@"
$c.Value = '29144'
$c.Style = "Normal"
$c = $ws2.Cells.Item(4,8)
$c.NumberFormat = "@"
$c.Value = '119933'
$c.Style = "Normal"
$c = $ws2.Cells.Item(4,9)
$c.NumberFormat = "@"
$c.Value = '260194944'
$c.Style = "Normal"
$c = $ws2.Cells.Item(4,10)
$c.NumberFormat = "@"
$c.Value = '03/17/2015'
$c.Style = "Normal"
$c = $ws2.Cells.Item(4,11)
$c.NumberFormat = "@"
$c.Value = 'The college student with Epilepsy who got ripped off!'
$c.Style = "Normal"
$c = $ws2.Cells.Item(4,12)
$c.NumberFormat = "@"
$c.Value = 'My spouse and I decided to go on a couples retreat and this location not only took our money but advised us to call a customer concerns number which has never returned any of our phone calls to even to pretend they care. we never received an apology from Jessie the manager but he will tell you to leave his property with police threats just to keep him from doing his job! I''m not a criminal but he stole from me and spouse. The whole stay we had a shower which never even got warm and it and never offered us another room but had put in maintenance request like I was living there.Maintenance stressed how  he couldn''t get over time so he couldn''t  yourour room shower today! If your smart take heed please as this place will take your deposit and just make you feel low and I know I''m blessed as a child of God. I work hard and am not rich but doesn''t mean I don''t deserve to receive what I pay for!MoreShow less'
$c.Style = "Normal"
$ws2.Cells.Item(4,13).Value = 1
$c = $ws2.Cells.Item(4,14)
$c.NumberFormat = "@"
$c.Value = 'March 2015'
$c.Style = "Normal"
$c = $ws2.Cells.Item(4,15)
$c.NumberFormat = "@"
$c.Value = ' traveled with family'
$c.Style = "Normal"
$ws2.Cells.Item(4,22).Value = 0
$c = $ws2.Cells.Item(4,25)
$c.NumberFormat = "@"
$c.Value = 'My spouse and I decided to go on a couples retreat and this location not only took our money but advised us to call a customer concerns number which has never returned any of our phone calls to even to pretend they care. we never received an apology from Jessie the manager but he will tell you to leave his property with police threats just to keep him from doing his job! I''m not a criminal but he stole from me and spouse. The whole stay we had a shower which never even got warm and it and never offered us another room but had put in maintenance request like I was living there.Maintenance stressed how  he couldn''t get over time so he couldn''t  yourour room shower today! If your smart take heed please as this place will take your deposit and just make you feel low and I know I''m blessed as a child of God. I work hard and am not rich but doesn''t mean I don''t deserve to receive what I pay for!More'
$c.Style = "Normal"

# row 5
$ws2.Cells.Item(5,1).Value = 39431
$ws2.Cells.Item(5,4).Value = 4
$c = $ws2.Cells.Item(5,5)
$c.NumberFormat = "@"
$c.Value = '08/03/2018'
$c.Style = "Normal"
$c = $ws2.Cells.Item(5,6)
$c.NumberFormat = "@"
$c.Value = 'https://www.tripadvisor.com/ShowUserReviews-g29144-d119933-r242853273-InTown_Suites_Denver_East-Aurora_Colorado.html'
$c.Style = "Normal"
$c = $ws2.Cells.Item(5,7)
$c.NumberFormat = "@"
$c.Value = '29144'
$c.Style = "Normal"
$c = $ws2.Cells.Item(5,8)
$c.NumberFormat = "@"
$c.Value = '119933'
$c.Style = "Normal"
$c = $ws2.Cells.Item(5,9)
$c.NumberFormat = "@"
$c.Value = '242853273'
$c.Style = "Normal"
$c = $ws2.Cells.Item(5,10)
$c.NumberFormat = "@"
$c.Value = '12/03/2014'
$c.Style = "Normal"
$c = $ws2.Cells.Item(5,11)
$c.NumberFormat = "@"
$c.Value = 'DO NOT STAY HERE!!'
$c.Style = "Normal"
$c = $ws2.Cells.Item(5,12)
$c.NumberFormat = "@"
$c.Value = 'I tried to give this place zero stars, but I had to pick at least one. Old dirty in need of repair. Never went without shoes in my own room. Even to the bathroom. I never put my luggage on the ground. I opened it on the table and cabinet. I wouldn''t take anything out unless I was putting it on my body. There is no dishes or pots and pans or utensils. Towels are paper thin and small. Maid came today and didn''t even make the bed!  NO BLANKET. Yes there was a comforter. Had to go buy a blanket at target. Will not refund my money to get out early. Dogs barking, kids crying or crawling up and down the hall. Smells coming out of rooms make me want to puke. I feel like the bed is going to fall it is so shaky. The pillows are plastic like at a hospital. Furniture is old and falling apart. Carpet needs to be thrown away. This is not a budget hotel. This is a I have no job no money no car to sleep in hotel. I can''t wait to move out of here. They will not give a refund so if I leave I"m out $450 and then go pay to stay somewhere else. I can''t afford that. Save yourself and your money and RUN AWAY FROM THIS PROPERTY!!!
In the time it took me to right this...I tried to give this place zero stars, but I had to pick at least one. Old dirty in need of repair. Never went without shoes in my own room. Even to the bathroom. I never put my luggage on the ground. I opened it on the table and cabinet. I wouldn''t take anything out unless I was putting it on my body. There is no dishes or pots and pans or utensils. Towels are paper thin and small. Maid came today and didn''t even make the bed!  NO BLANKET. Yes there was a comforter. Had to go buy a blanket at target. Will not refund my money to get out early. Dogs barking, kids crying or crawling up and down the hall. Smells coming out of rooms make me want to puke. I feel like the bed is going to fall it is so shaky. The pillows are plastic like at a hospital. Furniture is old and falling apart. Carpet needs to be thrown away. This is not a budget hotel. This is a I have no job no money no car to sleep in hotel. I can''t wait to move out of here. They will not give a refund so if I leave I"m out $450 and then go pay to stay somewhere else. I can''t afford that. Save yourself and your money and RUN AWAY FROM THIS PROPERTY!!!In the time it took me to right this review and get it to load I had time to color my hair, take a shower and warm up dinner. The shower has such low water pressure It takes twice as long to take a shower since I have thick hair. So not much of a energy saver there. I hope you like either a cold or hot shower. There is no in between. Keep driving……away from this place!MoreShow less'
$c.Style = "Normal"
$ws2.Cells.Item(5,13).Value = 1
$c = $ws2.Cells.Item(5,14)
$c.NumberFormat = "@"
$c.Value = 'December 2014'
$c.Style = "Normal"
$c = $ws2.Cells.Item(5,15)
$c.NumberFormat = "@"
$c.Value = ' traveled on business'
$c.Style = "Normal"
$ws2.Cells.Item(5,18).Value = 1
$ws2.Cells.Item(5,21).Value = 1
$ws2.Cells.Item(5,22).Value = 0
$c = $ws2.Cells.Item(5,25)
$c.NumberFormat = "@"
$c.Value = 'I tried to give this place zero stars, but I had to pick at least one. Old dirty in need of repair. Never went without shoes in my own room. Even to the bathroom. I never put my luggage on the ground. I opened it on the table and cabinet. I wouldn''t take anything out unless I was putting it on my body. There is no dishes or pots and pans or utensils. Towels are paper thin and small. Maid came today and didn''t even make the bed!  NO BLANKET. Yes there was a comforter. Had to go buy a blanket at target. Will not refund my money to get out early. Dogs barking, kids crying or crawling up and down the hall. Smells coming out of rooms make me want to puke. I feel like the bed is going to fall it is so shaky. The pillows are plastic like at a hospital. Furniture is old and falling apart. Carpet needs to be thrown away. This is not a budget hotel. This is a I have no job no money no car to sleep in hotel. I can''t wait to move out of here. They will not give a refund so if I leave I"m out $450 and then go pay to stay somewhere else. I can''t afford that. Save yourself and your money and RUN AWAY FROM THIS PROPERTY!!!
In the time it took me to right this...I tried to give this place zero stars, but I had to pick at least one. Old dirty in need of repair. Never went without shoes in my own room. Even to the bathroom. I never put my luggage on the ground. I opened it on the table and cabinet. I wouldn''t take anything out unless I was putting it on my body. There is no dishes or pots and pans or utensils. Towels are paper thin and small. Maid came today and didn''t even make the bed!  NO BLANKET. Yes there was a comforter. Had to go buy a blanket at target. Will not refund my money to get out early. Dogs barking, kids crying or crawling up and down the hall. Smells coming out of rooms make me want to puke. I feel like the bed is going to fall it is so shaky. The pillows are plastic like at a hospital. Furniture is old and falling apart. Carpet needs to be thrown away. This is not a budget hotel. This is a I have no job no money no car to sleep in hotel. I can''t wait to move out of here. They will not give a refund so if I leave I"m out $450 and then go pay to stay somewhere else. I can''t afford that. Save yourself and your money and RUN AWAY FROM THIS PROPERTY!!!In the time it took me to right this review and get it to load I had time to color my hair, take a shower and warm up dinner. The shower has such low water pressure It takes twice as long to take a shower since I have thick hair. So not much of a energy saver there. I hope you like either a cold or hot shower. There is no in between. Keep driving……away from this place!More'
$c.Style = "Normal"
